# Convert the opening Heading1/bold "By ..." block into a pandoc-style
# title block: a Title-styled paragraph followed by an Authors-styled
# paragraph, each run split word-by-word (matching pandoc's output),
# and drop the bookmark that used to wrap the heading.

$d = $word.ActiveDocument

# --- Remove the bookmark that wrapped the old Heading1 paragraph ---
if ($d.Bookmarks.Count -gt 0) {
    $d.Bookmarks.Item(1).Delete()
}

# --- Paragraph 1: title ---
$p1 = $d.Paragraphs.Item(1)
$p1.Range.Text = "Distributism Versus Capitalism"
$p1.Style = "Title"

# Split the title text into separate runs: "Distributism" " " "Versus" " " "Capitalism"
$p1Range = $p1.Range
$p1Start = $p1Range.Start
$p1Range.Text = ""
$d.Range($p1Start, $p1Start).InsertAfter("Distributism")
$d.Range($p1Start + 12, $p1Start + 12).InsertAfter(" ")
$d.Range($p1Start + 13, $p1Start + 13).InsertAfter("Versus")
$d.Range($p1Start + 19, $p1Start + 19).InsertAfter(" ")
$d.Range($p1Start + 20, $p1Start + 20).InsertAfter("Capitalism")

# --- Paragraph 2: authors ---
$p2 = $d.Paragraphs.Item(2)
$p2.Range.Text = "Dorothy Day"
$p2.Style = "Authors"
$p2.Range.Font.Bold = 0

$p2Range = $p2.Range
$p2Start = $p2Range.Start
$p2Range.Text = ""
$d.Range($p2Start, $p2Start).InsertAfter("Dorothy")
$d.Range($p2Start + 7, $p2Start + 7).InsertAfter(" ")
$d.Range($p2Start + 8, $p2Start + 8).InsertAfter("Day")

"done"
